$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Regions")

# Insert a new header row at the top of the data (row 3), pushing the
# existing state/region rows down by one.
$ws.Rows("3:3").Insert()

# Populate the new header row with column labels and make them bold.
$ws.Range("A3").Value = "state"
$ws.Range("B3").Value = "region"
$ws.Range("A3:B3").Font.Bold = $true

# Make "Regions" the active/visible sheet and select B9, matching the
# saved view state of the workbook.
$ws.Activate()
$ws.Range("B9").Select() | Out-Null
